$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fridge")

# --- Convert the StartDate / ExpireDate columns from text dates to real
#     date values (serial numbers), and correct a few expiration dates. ---
$ws.Range("B2").Value2 = 45658          # 01/01/2025
$ws.Range("C2").Value2 = 45667          # 10/01/2025

$ws.Range("B3").Value2 = 45664          # 07/01/2025
$ws.Range("C3").Value2 = 45699          # 11/02/2025

$ws.Range("B4").Value2 = 45672          # 15/01/2025
$ws.Range("C4").Value2 = 45700          # 12/02/2025

$ws.Range("B5").Value2 = 45668          # 11/01/2025
$ws.Range("C5").Value2 = 45686          # 29/01/2025

$ws.Range("B2:C5").NumberFormat = "dd/mm/yy;@"

# --- Rename status label "Expired Soon" -> "Expires Soon" ---
$ws.Range("D5").Value = "Expires Soon"

# --- Best-fit the StartDate column width to its new date content ---
$ws.Columns.Item(2).AutoFit()

# --- Page setup (paper size / orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection shown in the sheet view ---
$ws.Range("F10").Select()
